$d = $word.ActiveDocument
$d.Content.Find.Execute("meilleur performance", $false, $false, $false, $false, $false, $true, 1, $false, "meilleure performance", 2)
